$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the district names (shared strings)
$ws.Range("A1").Value = "Кировский"
$ws.Range("A2").Value = "Емельяновский"

# Update the numeric values
$ws.Range("B1").Value = 7
$ws.Range("B2").Value = 37

# Move the active selection to A2 (matches the saved selection in the sheet view)
$ws.Range("A2").Select()
